$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.913.70'
$ws.Range("E2").Value = '  +0.20%  '

$ws.Range("D3").Value = '3.788.83'
$ws.Range("E3").Value = '  -0.96%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.99'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.22'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.14%  '

$ws.Range("D7").Value = '3.786.63'
$ws.Range("E7").Value = '  -1.03%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.514'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.157'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.84'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +8.54%  '

$ws.Range("E12").Value = '  -1.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000246'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.02'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.55%  '

$ws.Range("D15").Value = '4.427.65'
$ws.Range("E15").Value = '  -0.95%  '

$ws.Range("D16").Value = '3.793.60'
$ws.Range("E16").Value = '  -0.59%  '

$ws.Range("D17").Value = '67.915.15'
$ws.Range("E17").Value = '  +0.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.17'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.01'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '458.26'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.45'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.690'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.60%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.17'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.32%  '

$ws.Range("E25").Value = '  -4.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.87'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.85%  '

$ws.Range("E27").Value = '  -1.90%  '

$ws.Range("E28").Value = '  -0.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.90'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.10%  '

$ws.Range("D30").Value = '3.940.17'
$ws.Range("E30").Value = '  -0.84%  '

$ws.Range("E31").Value = '  -7.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.20'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.76%  '

$ws.Range("E33").Value = '  -1.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.98'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.47%  '

$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.92'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.83%  '

$ws.Range("E37").Value = '  -1.00%  '

$ws.Range("E38").Value = '  +5.81%  '

$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("E40").Value = '  -4.39%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.980'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.09%  '

$ws.Range("E43").Value = '  +0.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.63'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.13'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '152.29'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.74%  '

$ws.Range("E47").Value = '  -2.34%  '

$ws.Range("E48").Value = '  -2.62%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.28'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.99%  '

$ws.Range("E50").Value = '  -1.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.38'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -8.54%  '

